# Empower-Remote-Workers-Poster.pptx edit
#
# 1) The cached "datetimeFigureOut" field text (date last refreshed by
#    PowerPoint) moves from 10/23/2019 to 10/30/2019 everywhere it is
#    cached: the slide master, all 11 slide layouts, and the notes master.
# 2) A stray double space in the "Azure  AD" caption on the poster slide
#    is collapsed to a single space.

$p = $ppt.ActivePresentation

$oldDate = "10/23/2019"
$newDate = "10/30/2019"

function Update-DateField($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# Slide master's cached date field.
$master = $p.SlideMaster
Update-DateField $master.Shapes

# Every slide layout has its own cached copy of the date field.
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    Update-DateField $layout.Shapes
}

# Notes master's cached date field.
$notesMaster = $p.NotesMaster
Update-DateField $notesMaster.Shapes

# Fix the doubled space in the "Azure  AD" label on the poster slide.
$slide = $p.Slides.Item(1)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "Azure  AD") {
            $tr.Text = "Azure AD"
        }
    }
}
